# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-22
$newValues = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 3
    20 = 3
    21 = 3
    22 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
